# [#5] Begin sheets reinitialization on deconstruction
#
# - Rename Sheet1 -> Timesheet, Sheet2 -> Cfg
# - Remove Sheet3
# - Merge A1:B1 on Timesheet, left-align the merged cell
# - Freeze panes on Timesheet at C3 (2 cols / 2 rows), restore per-pane
#   selections (topRight=C1, bottomLeft=A3, bottomRight=C10)
# - Set selection on Cfg to H7

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")

$ws1.Name = "Timesheet"
$ws2.Name = "Cfg"

$excel.DisplayAlerts = $false
$ws3.Delete() | Out-Null
$excel.DisplayAlerts = $true

# Left-align A1:B1 first, then merge, so both cells end up sharing the
# same (new) style index.
$ws1.Range("A1:B1").HorizontalAlignment = -4131
$ws1.Range("A1:B1").Merge() | Out-Null

# Set up the Cfg sheet's selection before switching back to Timesheet so
# that Timesheet ends up as the active/selected tab.
$ws2.Activate()
$ws2.Range("H7").Select() | Out-Null

# Freeze panes on Timesheet with top-left of the scrolling area at C3,
# then move the active selection to C10.
$ws1.Activate()
$ws1.Range("C3").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws1.Range("C10").Select() | Out-Null
